$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip back the per-category "moc" band (L,M,O:S,U:X) and "day" (Z:AC)
# detail columns for the existing data rows (2-5). Only the summary
# columns K (mocSMS), N (mocRobo), T (bandT6), Y (dayTue), AD (daySun)
# and AE (datNat) survive the de-scrub.
$ws.Range("J2:J5").ClearContents()
$ws.Range("L2:M5").ClearContents()
$ws.Range("O2:S5").ClearContents()
$ws.Range("U2:X5").ClearContents()
$ws.Range("Z2:AC5").ClearContents()

# Append the newly-integrated record as row 6, using the same reduced
# column layout as the scrubbed rows above it.
$a6 = $ws.Range("A6")
$a6.Value = "'9123456809"
$a6.Style = "Normal"

$ws.Range("K6").Value = $true
$ws.Range("N6").Value = $true
$ws.Range("T6").Value = $true
$ws.Range("Y6").Value = $true
$ws.Range("AD6").Value = $true
$ws.Range("AE6").Value = $true
